# Collapse the three identical empty "Calibri" paragraphs (right after the
# "<<daysPlus20>>" paragraph) down to a single paragraph, and move the
# "_GoBack" bookmark from the final paragraph of the document onto that
# surviving empty paragraph.

$d = $word.ActiveDocument

# The three empty paragraphs are items 6, 7, 8 in the Paragraphs collection.
# Delete the 3rd and 2nd (highest index first so earlier indices stay valid),
# leaving a single paragraph (item 6) in their place.
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# Moving a bookmark named "_GoBack" automatically removes it from wherever it
# previously lived (the last paragraph before the sectPr), matching Word's
# single-instance behavior for this bookmark.
$d.Bookmarks.Add("_GoBack", $d.Paragraphs.Item(6).Range)
